$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (2..302).
# The workbook was refreshed one day later, so bump the date from
# 2023-10-03 (45202) to 2023-10-04 (45203) for all of them.
$ws.Range("C2:C302").Value = 45203
